# [Kadastro App] Yeni kayit eklendi: 2937
# Adds a new record row (2937 / 2025-09-09 / Erdemli / 1 / CAP / SEVIL SARACER)
# to both the "Kayitlar" master sheet and the "Erdemli" district sheet.

$wb = $excel.ActiveWorkbook

function Add-KayitRow {
    param($ws, $row, $kayitNo, $tarih, $birim, $parselSayisi, $isAdi, $personeller)

    $values = @($kayitNo, $tarih, $birim, $parselSayisi, $isAdi, $personeller)
    $cols = @("A", "B", "C", "D", "E", "F")

    for ($i = 0; $i -lt $cols.Length; $i++) {
        $cell = $ws.Range($cols[$i] + "$row")
        # Build the literal value via a text formula so Excel never
        # auto-coerces numeric-looking text (e.g. "2937", "1") or date-
        # looking text (e.g. "2025-09-09") into a number/date ...
        $escaped = $values[$i].Replace('"', '""')
        $cell.Formula = '="' + $escaped + '"'
        # ... then collapse the formula down to its plain string result,
        # so the stored cell is a bare literal (no <f>, no extra style).
        $cell.Copy()
        $cell.PasteSpecial(-4163)
    }
}

$wsKayitlar = $wb.Worksheets.Item("Kayitlar")
$kayitlarNextRow = $wsKayitlar.UsedRange.Rows.Count + 1
Add-KayitRow $wsKayitlar $kayitlarNextRow "2937" "2025-09-09" "Erdemli" "1" "ÇAP" "SEVİL SARAÇER (Tekniker)"

$wsErdemli = $wb.Worksheets.Item("Erdemli")
$erdemliNextRow = $wsErdemli.UsedRange.Rows.Count + 1
Add-KayitRow $wsErdemli $erdemliNextRow "2937" "2025-09-09" "Erdemli" "1" "ÇAP" "SEVİL SARAÇER (Tekniker)"
